$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.046.86"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "2.921.07"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'199.46"
$ws.Range("E5").Value = "  +5.74%  "
$ws.Range("D6").Value = "'599.22"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("D9").Value = "'0.198"
$ws.Range("E9").Value = "  +2.08%  "
$ws.Range("D10").Value = "2.917.61"
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("E11").Value = "  +15.69%  "
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "3.455.44"
$ws.Range("E14").Value = "  +2.56%  "
$ws.Range("D15").Value = "75.906.97"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("E16").Value = "  +1.81%  "
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "2.916.85"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("E19").Value = "  -2.22%  "
$ws.Range("D20").Value = "'12.79"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("D21").Value = "'378.26"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("D24").Value = "'71.43"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "3.067.96"
$ws.Range("E26").Value = "  +2.83%  "
$ws.Range("D27").Value = "'4.23"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "'9.73"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("E29").Value = "  +5.34%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").Value = "'506.52"
$ws.Range("E32").Value = "  -4.69%  "
$ws.Range("D33").Value = "'7.73"
$ws.Range("E33").Value = "  -2.89%  "
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'164.88"
$ws.Range("D37").Value = "'20.17"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").Value = "'19.70"
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("E39").Value = "  -6.03%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'180.14"
$ws.Range("E41").Value = "  -2.51%  "
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.0965"
$ws.Range("E43").Value = "  +11.82%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").Value = "'5.00"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("D46").Value = "'40.19"
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").Value = "'1.20"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").Value = "'0.577"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  +7.27%  "
$ws.Range("E51").Value = "  -1.14%  "
